# Update "想去人数" (column F) values on both the "展览" sheet and the
# "全部类型" sheet, keyed by the event name in column C so the correct
# row is updated on each sheet even though row numbers differ slightly
# between the two sheets.

$wb = $excel.ActiveWorkbook

# name (column C) -> new value for column F
$updates = @{
    "杭州·SST动漫嘉年华" = 1586
    "杭州·AD02动漫展" = 9635
    "杭州·星玫Rostar偶像团 1st off会 - 莫里生日SP" = 15
    "杭州·幻想物语动漫游戏展" = 1205
    "杭州·2024ESCC游戏电竞博览会暨新次元微光青春动漫交流会" = 2998
    "杭州·ELECTRIC COMIC动漫游戏展" = 2263
    "杭州·梦漫星河动漫展" = 1962
    "杭州·第36届 中二病 原神x星穹only" = 1564
    "杭州·ACG发色only-黑白两色" = 309
    "杭州·young girls二次元全女夜场" = 21
    "杭州·代号鸢相聚广陵2.0only（中婚版）" = 187
    "杭州·赛马娘only—晴空雏菊" = 221
    "杭州·SK怀旧展&偶像专场" = 346
    "杭州·白日梦次元动漫嘉年华" = 316
    "杭州·造梦探险家——二次元同好会" = 34
    "杭州·与梦回望动漫游戏展" = 147
    "杭州·第37届 中二病 原神x星穹only" = 1543
    "杭州·AD03动漫展" = 175
    "杭州·第四届华盟动漫次元嘉年华" = 1511
    "杭州·造梦探险家Porject6野蛮冲撞——第五人格ONLY" = 41
    "杭州·现世繁华-代号鸢only" = 355
    "杭州·第八届YH樱花动漫游戏文化节" = 381
    "杭州·第三届日夜国乙only" = 781
}

# "杭州·Look Look动漫嘉年华" appears twice (two separate rows with the
# same name) on both sheets, both changing from their old value to 700.
$lookLookValue = 700

foreach ($ws in $wb.Worksheets) {
    $usedRange = $ws.UsedRange
    $lastRow = $usedRange.Rows.Count

    for ($r = 2; $r -le $lastRow; $r++) {
        $nameCell = $ws.Cells.Item($r, 3)   # column C
        $name = $nameCell.Value()

        if ($null -eq $name) { continue }

        if ($name -eq "杭州·Look Look动漫嘉年华") {
            $ws.Cells.Item($r, 6).Value = $lookLookValue
        }
        elseif ($updates.ContainsKey($name)) {
            $ws.Cells.Item($r, 6).Value = $updates[$name]
        }
    }
}
